$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 3, shifting existing rows 3..37 down to 4..38
$ws.Rows.Item(3).Insert()

# Populate the newly inserted row 3 with data.
$ws.Cells.Item(3, 1).Value = 4
$ws.Cells.Item(3, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(3, 3).Value = "Los Lagos"
$ws.Cells.Item(3, 4).Value = 44921
$ws.Cells.Item(3, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(3, 5).Value = 10
$ws.Cells.Item(3, 6).Value = "Fruta"
$ws.Cells.Item(3, 7).Value = 100101
$ws.Cells.Item(3, 8).Value = "Berries"
$ws.Cells.Item(3, 9).Value = 100101001
$ws.Cells.Item(3, 10).Value = "Arándano (blue)"
$ws.Cells.Item(3, 11).Value = "Sin especificar"
$ws.Cells.Item(3, 12).Value = "Primera"
$ws.Cells.Item(3, 13).Value = 200
$ws.Cells.Item(3, 14).Value = 3000
$ws.Cells.Item(3, 15).Value = 3500
$ws.Cells.Item(3, 16).Value = 3250
$ws.Cells.Item(3, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(3, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(3, 19).Value = 1625
$ws.Cells.Item(3, 20).Value = 2
